$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the window tab-ratio (cosmetic bookViews/workbookView@tabRatio, 176 -> 611 in OOXML
# thousandths, i.e. ~0.611 as a fraction on the COM Window object).
$win = $excel.ActiveWindow
$win.TabRatio = 0.611

# Refresh / get fixed QPSK (column C) values
$ws.Range("C2").Value = 0.0254
$ws.Range("C3").Value = 0.0183
$ws.Range("C4").Value = 0.0118
$ws.Range("C5").Value = 0.006
$ws.Range("C6").Value = 0.0022
$ws.Range("C7").Value = 0.0005
$ws.Range("C8").Formula = "=4/16000000"
$ws.Range("C9").Value = 0

# Move the selection/active cell to C9
$ws.Range("C9").Select()
